$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (A28:C28 - date
# style + wrapped-text note style + default type style) onto the two new
# rows so the new cells reuse the workbook's existing styles.
$ws.Range("A28:C28").Copy() | Out-Null
$ws.Range("A29:C29").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:C30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 29: new journal entry dated 2018-03-01 (serial 43160)
$ws.Cells.Item(29, 1).Value = 43160
$ws.Cells.Item(29, 2).Value = "Rédaction de la documentation du projet "
$ws.Cells.Item(29, 3).Value = "2 périodes"

# Row 30: second new journal entry, same day
$ws.Cells.Item(30, 1).Value = 43160
$ws.Cells.Item(30, 2).Value = "J'ai commencé la fonction qui me permettra d'ajouter des nouveaux articles dans la base de données "
$ws.Cells.Item(30, 3).Value = "2 périodes"

# Row 30's note is long enough to wrap onto two lines, so its row grows
# to double height like the other wrapped-text rows in the sheet.
$ws.Rows.Item(30).RowHeight = 30

# Move the selection past the newly-added rows.
[void]$ws.Range("C31").Select()
